$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Nganh" (major) column D values: KTPM / ATTT -> "Công nghệ phần mềm"; Marketing -> "Ngành giáo dục mầm non"
$ws.Range("D2").Value = "Công nghệ phần mềm"
$ws.Range("D3").Value = "Công nghệ phần mềm"
$ws.Range("D5").Value = "Công nghệ phần mềm"

$ws.Range("D4").Value = "Ngành giáo dục mầm non"
$ws.Range("D6").Value = "Ngành giáo dục mầm non"
$ws.Range("D7").Value = "Ngành giáo dục mầm non"

# Update "Khoa" (department) column C values: KT -> GDMN
$ws.Range("C4").Value = "GDMN"
$ws.Range("C6").Value = "GDMN"
$ws.Range("C7").Value = "GDMN"

# Widen column D to fit the new, longer major names
$ws.Columns.Item(4).ColumnWidth = 38.33

# Update the active selection to I6 (matches author's last cursor position)
$ws.Range("I6").Select()
